# "Taking back Chinna's changes" — revert the Runmode/Result columns on the
# "Test Cases" sheet back to "N" / real results, and add the
# AuthoringRecordViewDetailsTest row that had been dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 2: Result goes from SKIP -> FAIL
$ws.Range("D2").Value = "FAIL"

# Rows 3-5: Runmode goes from Y -> N
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"

# New row 6: AuthoringRecordViewDetailsTest test case.
# Copy the formatting from existing cells first so the new row matches
# the look of the rest of the table, then fill in the values.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null

$ws.Range("D5").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null

$ws.Range("C4").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null

$ws.Range("D5").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null

$ws.Range("A6").Value = "AuthoringRecordViewDetailsTest"
$ws.Range("B6").Value = "To verify Record View Details link Navigate to WOS page and Navigate to Project Neon Page"
$ws.Range("C6").Value = "N"
$ws.Range("D6").Value = "SKIP"

# Restore the selection to what it was (C12) on this sheet.
$ws.Range("C12").Select()
